$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.734148025512695
$ws.Range("B1").Value = 2.263491630554199
$ws.Range("C1").Value = 1.231330513954163
$ws.Range("D1").Value = 1.304223775863647
$ws.Range("E1").Value = 1.50303316116333
